# Insert a new weekly price record as row 28, pushing all subsequent
# records down by one row (dimension grows from A1:R129 to A1:R130).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 28..129 down to 29..130.
$ws.Rows.Item(28).Insert()

# Populate the newly-inserted row 28 with the new record.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44742
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 100112052
$ws.Range("G28").Value = "Albahaca"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 7500
$ws.Range("N28").Value = '$/paquete'
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 7500
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
